$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in F1, matching the style of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the time_taken column for each data row (F2:F16) with the recorded timestamps
$ws.Range("F2").Value = "2021-10-05 10:50:40.347558"
$ws.Range("F3").Value = "2021-10-05 10:50:40.347570"
$ws.Range("F4").Value = "2021-10-05 10:50:40.347574"
$ws.Range("F5").Value = "2021-10-05 10:50:40.347577"
$ws.Range("F6").Value = "2021-10-05 10:50:40.347581"
$ws.Range("F7").Value = "2021-10-05 10:50:40.347584"
$ws.Range("F8").Value = "2021-10-05 10:50:40.347587"
$ws.Range("F9").Value = "2021-10-05 10:50:40.347590"
$ws.Range("F10").Value = "2021-10-05 10:50:40.347593"
$ws.Range("F11").Value = "2021-10-05 10:50:40.347597"
$ws.Range("F12").Value = "2021-10-05 10:50:40.347600"
$ws.Range("F13").Value = "2021-10-05 10:50:40.347603"
$ws.Range("F14").Value = "2021-10-05 10:50:40.347606"
$ws.Range("F15").Value = "2021-10-05 10:50:40.347609"
$ws.Range("F16").Value = "2021-10-05 10:50:40.347612"
